# Insert a "-------------" separator row right above the
# "PORCENT_COMPLETA (%)" / "PORCENT_AUSENTE (%)" summary rows on the first
# four worksheets (Planilha1..Planilha4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Planilha1", "Planilha2", "Planilha3", "Planilha4")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Find the row that currently holds "PORCENT_COMPLETA (%)" in column A.
    $used = $ws.UsedRange
    $lastRow = $used.Rows.Count
    $targetRow = -1

    for ($r = 1; $r -le $lastRow; $r++) {
        $val = $ws.Cells.Item($r, 1).Value()
        if ($val -eq "PORCENT_COMPLETA (%)") {
            $targetRow = $r
            break
        }
    }

    if ($targetRow -gt 0) {
        # Push PORCENT_COMPLETA / PORCENT_AUSENTE rows down by one, then
        # write the separator label into the freed row.
        $ws.Rows.Item($targetRow).Insert()
        $cell = $ws.Cells.Item($targetRow, 1)
        $cell.Style = "Normal"
        $cell.Value = "-------------"
    }
}
